$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Explorers
$ws.Range("B2").Value = 114
$ws.Range("D2").Value = 6.06
$ws.Range("E2").Value = 5.34
$ws.Range("F2").Value = 2.09
$ws.Range("G2").Value = 2.62
$ws.Range("H2").Value = 5.67

# Row 3 - Pioneers
$ws.Range("B3").Value = 226
$ws.Range("C3").Value = 0.26
$ws.Range("D3").Value = 5.2
$ws.Range("E3").Value = 4.84
$ws.Range("G3").Value = 3.85
$ws.Range("H3").Value = 4.79

# Row 4 - Hesitators
$ws.Range("B4").Value = 352
$ws.Range("D4").Value = 4.62
$ws.Range("E4").Value = 3.05
$ws.Range("F4").Value = 3.65
$ws.Range("G4").Value = 3.9
$ws.Range("H4").Value = 4.03

# Row 5 - Avoiders
$ws.Range("B5").Value = 188
$ws.Range("C5").Value = 0.21
$ws.Range("D5").Value = 3.8
$ws.Range("E5").Value = 1.58
$ws.Range("F5").Value = 4.76
$ws.Range("H5").Value = 2.98
